# ajout incertitude sur la pente
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the slope value B2 (previously the formula =6*10^-4) with the new
# literal uncertainty-adjusted value, removing the formula.
$ws.Range("B2").Value = 0.0011000000000000001

# The remaining B2:B5 cells held simple literal formulas (e.g. =1.5*10^-3);
# convert them to plain values, keeping the same numeric results.
$ws.Range("B3").Value = 0.0015
$ws.Range("B4").Value = 0.002
$ws.Range("B5").Value = 0.0027000000000000001

# Column C held a shared formula computing the resistance; replace every
# cell with its last calculated value so no formula remains.
$ws.Range("C2").Value = 0.00055016523537939125
$ws.Range("C3").Value = 0.0011003304707587825
$ws.Range("C4").Value = 0.0016504957061381736
$ws.Range("C5").Value = 0.002200660941517565
$ws.Range("C6").Value = 0.0027508261768969559
$ws.Range("C7").Value = 0.0033009914122763473
$ws.Range("C8").Value = 0.003851156647655739

# Move the active selection to B3, matching the saved UI state.
$ws.Range("B3").Select()
